# Generate Report for Handback
# Updates the "Latest HO Xliff Generate Date" / handoff / handback timestamps
# for the c6f0ec6f-e1ce-4f66-ac63-826cb88a5aaf.md row (row 5) across the
# Overview, zh-cn and de-de worksheets, reflecting a new handback report run.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G5").Value = "2017-02-21 09:58:01"

# --- zh-cn sheet ------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H5").Value = "2017-02-21 09:57:43"
$wsZhCn.Range("L5").Value = "2017-02-21 09:58:44"

# --- de-de sheet ------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H5").Value = "2017-02-21 09:58:01"
$wsDeDe.Range("L5").Value = "2017-02-21 09:59:07"
